# Updated cryptos list on Sun Jul 21 03:45:45 UTC 2024 with GitHub Actions
# Refresh the coin price / volume figures (and fix the ONDO/Arweave row order)
# on the active worksheet of the already-open workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.238.21'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").Value = '3.513.98'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.23%  '
$ws.Range("E9").Value = '  +5.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '4.122.82'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.52%  '
$ws.Range("D15").Value = '67.169.98'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").Value = '3.531.30'
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '396.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.539'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("E34").Value = '  +4.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.41%  '
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0754'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.10%  '
$ws.Range("D43").Value = '2.835.70'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.54%  '
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").Value = '  -2.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '340.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.27%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.04%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("E51").Value = '  -0.31%  '
